$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7: 2021年 ---
$ws.Range("A7").Value = "2021年"
$ws.Range("B7").Value = 100.3
$ws.Range("C7").Value = 100.7
$ws.Range("D7").Value = 99.5

# Apply the same formatting used for the other year cells in column A
# (bold, bordered, centered) by copying formats from A6.
$ws.Range("A6").Copy()
$ws.Range("A7").PasteSpecial(-4122)  # xlPasteFormats

# --- Row 8: 2022年 ---
$ws.Range("A8").Value = "2022年"
$ws.Range("B8").Value = 100.6

$ws.Range("A6").Copy()
$ws.Range("A8").PasteSpecial(-4122)  # xlPasteFormats

# C8 and D8 are present but empty (no value). Force the cells to be
# materialized in the sheet (rather than left completely absent) by
# briefly toggling a border and removing it again.
$ws.Range("C8").Borders.LineStyle = 1
$ws.Range("C8").Borders.LineStyle = -4142  # xlLineStyleNone
$ws.Range("D8").Borders.LineStyle = 1
$ws.Range("D8").Borders.LineStyle = -4142  # xlLineStyleNone

$excel.CutCopyMode = $false
